# Revert "DOMA-2542 Localization for Excel template (ticket_report_status_executor)"
#
# The localization commit had rewritten `{d.tickets[i+1].<field>}` template
# placeholders (no spaces around the `+`) into `{d.tickets[i + 1].<field>}`
# (spaces around the `+`). This reverts that specific textual change on the
# worksheet's third row (A3:H3), which is the only row using the
# `tickets[i+1]` placeholder family.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "{d.tickets[i + 1].categoryClassifier}"
$ws.Range("B3").Value = "{d.tickets[i + 1].address}"
$ws.Range("C3").Value = "{d.tickets[i + 1].processing}"
$ws.Range("D3").Value = "{d.tickets[i + 1].completed}"
$ws.Range("E3").Value = "{d.tickets[i + 1].canceled}"
$ws.Range("F3").Value = "{d.tickets[i + 1].deferred}"
$ws.Range("G3").Value = "{d.tickets[i + 1].closed}"
$ws.Range("H3").Value = "{d.tickets[i + 1].new_or_reopened}"

# The same revert also restored the workbook theme's minor font from
# "+mn-lt"/"Helvetica Neue" back to the template's original "Cambria" (used
# by the default shape/text styles baked into the theme) and restored the
# drop-shadow effect style that had been cleared out. The COM surface here
# only exposes the theme's font scheme (not the effect styles), so update
# what is reachable: the workbook's minor theme font.
$wb.Theme.ThemeFontScheme.MinorFont.Latin = "Cambria"
